# Update localization-status report:
#  - Status text "Ready for handoff" -> "In Translation" on every sheet
#  - Narrow the "Status"/"zh-cn"/"de-de" columns that hold that text

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$targetWidth = 12.5   # closest achievable ColumnWidth producing the new, narrower column

# --- Overview sheet: status shown in columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $targetWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetWidth

# --- zh-cn sheet: status shown in column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $targetWidth

# --- de-de sheet: status shown in column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $targetWidth
